# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.022.03'
$ws.Range("E2").Value = '  +2.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.939.89'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.43'
$ws.Range("E5").Value = '  +3.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.13'
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("E7").Value = '  -1.41%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.724'
$ws.Range("E9").Value = '  -2.36%  '

$ws.Range("E10").Value = '  +8.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000356'
$ws.Range("E11").Value = '  +13.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.72'
$ws.Range("E12").Value = '  -1.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.573.52'
$ws.Range("E13").Value = '  +0.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.50'
$ws.Range("E14").Value = '  +1.27%  '

$ws.Range("B15").Value = 'Uniswap'
$ws.Range("C15").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.62'
$ws.Range("E15").Value = '  -1.15%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.940.53'
$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("E17").Value = '  -0.23%  '

$ws.Range("E18").Value = '  -1.76%  '

$ws.Range("E19").Value = '  -2.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.097.95'
$ws.Range("E20").Value = '  +2.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.88'
$ws.Range("E21").Value = '  +0.61%  '

$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.37'
$ws.Range("E22").Value = '  +2.83%  '

$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.60'
$ws.Range("E23").Value = '  -1.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.80'
$ws.Range("E24").Value = '  -0.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.67'
$ws.Range("E25").Value = '  +16.56%  '

$ws.Range("E26").Value = '  -0.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.36'
$ws.Range("E27").Value = '  +1.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.91'
$ws.Range("E28").Value = '  +8.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.20'
$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '711.41'
$ws.Range("E30").Value = '  -2.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.29'
$ws.Range("E31").Value = '  -3.06%  '

$ws.Range("E32").Value = '  -3.96%  '

$ws.Range("E33").Value = '  +3.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0928'
$ws.Range("E34").Value = '  +34.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.40'
$ws.Range("E35").Value = '  -4.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.75'
$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.151'
$ws.Range("E37").Value = '  -7.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.66'
$ws.Range("E38").Value = '  +3.25%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("E40").Value = '  -1.66%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.05'
$ws.Range("E41").Value = '  +9.75%  '

$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  +6.81%  '

$ws.Range("E43").Value = '  +3.01%  '

$ws.Range("E44").Value = '  -2.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.142'
$ws.Range("E45").Value = '  +0.39%  '

$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("E47").Value = '  -1.57%  '

$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '148.11'
$ws.Range("E49").Value = '  +2.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.13'
$ws.Range("E50").Value = '  -4.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.82'
$ws.Range("E51").Value = '  -1.54%  '
